$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.012.15'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.637.04'
$ws.Range('E3').Value = '  -0.37%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  -1.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.89%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.31'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0789'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '1.864.07'
$ws.Range('E12').Value = '  -0.39%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.659.22'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('E15').Value = '  -3.18%  '
$ws.Range('D16').Value = '25.994.85'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0₃0741'
$ws.Range('E17').Value = '  -3.13%  '
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.79%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('E23').Value = '  -2.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.132'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.79'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.74%  '
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('E30').Value = '  -1.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0481'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.04%  '
$ws.Range('E32').Value = '  -4.00%  '
$ws.Range('E33').Value = '  -4.75%  '
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.48'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.05%  '
$ws.Range('D36').Value = '1.135.56'
$ws.Range('E36').Value = '  +0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.863'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.65%  '
$ws.Range('E38').Value = '  -0.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.519'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.86%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0155'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.27'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.774'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '1.773.93'
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('E44').Value = '  -4.71%  '
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '54.95'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.51'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('E51').Value = '  +0.29%  '
